# Tripadvisor New Orleans shard 115 update:
#  1. Reorder the worksheet tabs so "review_info" comes before "hotel_info".
#  2. Add a new "State" column to "hotel_info" (inserted right after
#     "Hotel_Name", before "City") and populate it with "Louisiana" for
#     the existing hotel row.

$wb = $excel.ActiveWorkbook

# --- 1. Move "review_info" to be the first sheet tab ---------------------
$wsReview = $wb.Worksheets.Item("review_info")
$wsReview.Move($wb.Worksheets.Item(1))

# --- 2. Insert the "State" column into "hotel_info" -----------------------
$wsHotel = $wb.Worksheets.Item("hotel_info")

# Column C currently holds "City"; inserting here shifts City/Zip/etc. one
# column to the right, matching Hotel_Name | State | City | Zip | ...
$wsHotel.Range("C1").EntireColumn.Insert()

$wsHotel.Range("C1").Value = "State"
$wsHotel.Range("C2").Value = "Louisiana"
